$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.781.53"
$ws.Range("E2").Value = "  +3.37%  "
$ws.Range("D3").Value = "2.446.64"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.73"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.70"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "2.445.16"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").Value = "  +2.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.23"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +7.33%  "
$ws.Range("E15").Value = "  +5.80%  "
$ws.Range("D16").Value = "2.889.11"
$ws.Range("E16").Value = "  +3.12%  "
$ws.Range("D17").Value = "62.665.78"
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").Value = "2.450.25"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.89"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.99"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +2.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.45"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.03"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +8.81%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("B25").Value = "BabyDogeCoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D25").Value = "0.0₆0644"
$ws.Range("E25").Value = "  +131.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "66.01"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "652.86"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +12.72%  "
$ws.Range("E28").Value = "  +17.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.53"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +4.92%  "
$ws.Range("D30").Value = "0.0₃0988"
$ws.Range("E30").Value = "  +5.41%  "
$ws.Range("D31").Value = "2.567.68"
$ws.Range("E31").Value = "  +2.15%  "
$ws.Range("E32").Value = "  +8.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.21"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("E34").Value = "  +3.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.139"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +5.30%  "
$ws.Range("E36").Value = "  +2.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.77"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.52"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +6.68%  "
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "152.13"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  +2.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.71"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +8.00%  "
$ws.Range("E44").Value = "  +5.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.38"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  +27.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "145.43"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.72"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +6.70%  "
$ws.Range("E51").Value = "  +2.59%  "
